# Apply cryptocurrency price/volume updates to match the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# General format need to be forced to Text first, then the cell style is reset
# back to Normal/General so no visible formatting change is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.837.58"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.615.04"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue $ws.Range("D5") "212.34"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  +0.04%  "
Set-TextValue $ws.Range("D10") "19.77"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").Value = "1.841.46"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "1.608.21"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -0.73%  "
Set-TextValue $ws.Range("D15") "0.534"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "26.836.48"
$ws.Range("E16").Value = "  -1.10%  "
Set-TextValue $ws.Range("D17") "63.92"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.33%  "
Set-TextValue $ws.Range("D19") "210.56"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  +0.05%  "
Set-TextValue $ws.Range("D21") "6.76"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("E23").Value = "  -6.77%  "
Set-TextValue $ws.Range("D24") "8.94"
$ws.Range("E24").Value = "  -1.72%  "
Set-TextValue $ws.Range("D25") "146.40"
$ws.Range("E25").Value = "  -0.73%  "
Set-TextValue $ws.Range("D26") "7.50"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("E29").Value = "  -1.07%  "
Set-TextValue $ws.Range("D30") "0.0504"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("E32").Value = "  -2.54%  "
Set-TextValue $ws.Range("D33") "0.709"
$ws.Range("E33").Value = "  +30.00%  "
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").Value = "1.326.72"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  -1.29%  "
Set-TextValue $ws.Range("D37") "2.45"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -0.04%  "
Set-TextValue $ws.Range("D41") "0.795"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -0.54%  "
Set-TextValue $ws.Range("D44") "63.54"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "1.753.42"
$ws.Range("E45").Value = "  -1.34%  "
Set-TextValue $ws.Range("D46") "89.29"
$ws.Range("E46").Value = "  -1.52%  "
Set-TextValue $ws.Range("D47") "1.62"
$ws.Range("E47").Value = "  +1.12%  "
Set-TextValue $ws.Range("D48") "0.825"
$ws.Range("E48").Value = "  +8.12%  "
$ws.Range("E49").Value = "  -0.09%  "
Set-TextValue $ws.Range("D50") "0.0983"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.49"
$ws.Range("E51").Value = "  -1.58%  "
